$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.095.86"
$ws.Range("E2").Value = "  +4.37%  "

# Row 3
$ws.Range("D3").Value = "1.690.92"
$ws.Range("E3").Value = "  +3.25%  "

# Row 4
$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9966"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  -0.25%  "

# Row 5
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.11"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  +3.18%  "

# Row 6
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9975"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  -0.27%  "

# Row 7
$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4664"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  -1.08%  "

# Row 8
$ws.Range("E8").Value = "  +3.34%  "

# Row 9
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06177"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  +2.05%  "

# Row 10
$ws.Range("D10").Value = "1.681.28"
$ws.Range("E10").Value = "  +2.65%  "

# Row 11
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07033"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  +0.52%  "

# Row 12
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.19"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  +6.64%  "

# Row 13
$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.405"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  +2.97%  "

# Row 14
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5838"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  +3.55%  "

# Row 15
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "75.91"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  +3.98%  "

# Row 16
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9971"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  -0.31%  "

# Row 17
$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9978"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  -0.26%  "

# Row 18
$ws.Range("D18").Value = "26.085.23"
$ws.Range("E18").Value = "  +4.42%  "

# Row 19
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006752"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  +3.47%  "

# Row 20
$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.53"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  +3.14%  "

# Row 21
$ws.Range("D21").Value = "1.894.46"
$ws.Range("E21").Value = "  +2.44%  "

# Row 22
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.531"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  +6.91%  "

# Row 23
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.750"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  +3.89%  "

# Row 24
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.279"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  +2.09%  "

# Row 25
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "134.64"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  +1.71%  "

# Row 26
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.08"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  +2.21%  "

# Row 27
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.395"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  +2.70%  "

# Row 28
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.739"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  +7.22%  "

# Row 29
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "105.52"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  +2.40%  "

# Row 30
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.979"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  +2.88%  "

# Row 31
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.667"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  +5.28%  "

# Row 32
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07759"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  +2.99%  "

# Row 33
$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04371"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  +3.98%  "

# Row 34
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.594"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  +1.04%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6129"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  +4.78%  "

# Row 36
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9587"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  +3.51%  "

# Row 37
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9283"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  +5.92%  "

# Row 38
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "110.94"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  +14.00%  "

# Row 39
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.396"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  -6.45%  "

# Row 40
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9975"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -0.21%  "

# Row 41
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.889"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  +7.91%  "

# Row 42
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01463"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  -0.69%  "

# Row 43
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3752"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  +3.35%  "

# Row 44
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.059"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  +9.58%  "

# Row 45
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1132"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  +4.26%  "

# Row 46
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05316"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  +2.43%  "

# Row 47
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.204"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  +2.95%  "

# Row 48
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.79"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  +8.15%  "

# Row 49
$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.687"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  +8.86%  "

# Row 50
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.215"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  +3.32%  "

# Row 51
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9993"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  -0.20%  "
